# Update mods data [2025-12-26 15:10:42]
# Append a new daily record row (row 46) to the ModCounts sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ModCounts")

# Set the new row's values. The leading apostrophe on the date forces it
# to be stored as literal text (matching the existing Date column cells)
# instead of being auto-converted into a date serial number.
$ws.Range("A46").Value = "'2025/12/26"
$ws.Range("B46").Value = "逃离鸭科夫"
$ws.Range("C46").Value = 1102

# Match the formatting of the previous data row (center aligned style).
$ws.Range("A45:C45").Copy()
$ws.Range("A46:C46").PasteSpecial(-4122)
